# Atualização de bases das ligas, do dia: 29-03-2024 às 17:05
#
# The two most-recently-added matches (rows 86 and 87, i.e. id 84 and 85)
# had been entered with their data swapped: row 86 carried match
# 7511976 (DJK Bad Homburg vs SG Bornheim 1945 GrunWeiss) and row 87
# carried match 7511958 (SpVgg EGC Wirges vs SG 2000 MulheimKarlich),
# but the odds/result data for the two matches had been attached to the
# wrong id. This swaps all the match data (every column except the row
# id in column A) between the two rows so each id carries the correct
# match record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 86
$row2 = 87

$rangeA = $ws.Range("B$row1`:AC$row1")
$rangeB = $ws.Range("B$row2`:AC$row2")

$valuesA = $rangeA.Value2
$valuesB = $rangeB.Value2

$rangeA.Value2 = $valuesB
$rangeB.Value2 = $valuesA
